# Circle Language Spec Plan Part B, Sub-Projects.docx
# Reformulate into less resolute, more open language.

$d = $word.ActiveDocument
$WNS = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'
$CR = [char]13

function Find-ParagraphByText($text) {
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        $p = $d.Paragraphs($i)
        $t = $p.Range.Text.TrimEnd($CR)
        if ($t -eq $text) {
            return $p
        }
    }
    return $null
}

# ---------------------------------------------------------------------------
# 1. "Date: May 15, 2010 – May 15, 2010" -> "Date: May 15, 2010"
#    (drop the " – " run + duplicate smart-tagged date, reorder the
#    remaining smartTagPr attrs to Year/Day/Month)
# ---------------------------------------------------------------------------
$dateXml = @"
<w:p $WNS w:rsidR="002D5BF0" w:rsidRPr="00F012D7" w:rsidRDefault="00A37344">
  <w:pPr>
    <w:jc w:val="center"/>
    <w:rPr>
      <w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cs="Tahoma"/>
      <w:i/>
      <w:iCs/>
      <w:sz w:val="20"/>
      <w:szCs w:val="20"/>
      <w:lang w:val="en-US"/>
    </w:rPr>
  </w:pPr>
  <w:r w:rsidRPr="00F012D7">
    <w:rPr>
      <w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cs="Tahoma"/>
      <w:i/>
      <w:iCs/>
      <w:sz w:val="20"/>
      <w:szCs w:val="20"/>
      <w:lang w:val="en-US"/>
    </w:rPr>
    <w:t xml:space="preserve">Date: </w:t>
  </w:r>
  <w:smartTag w:uri="urn:schemas-microsoft-com:office:smarttags" w:element="date">
    <w:smartTagPr>
      <w:attr w:name="Year" w:val="2010"/>
      <w:attr w:name="Day" w:val="15"/>
      <w:attr w:name="Month" w:val="5"/>
    </w:smartTagPr>
    <w:r w:rsidRPr="00F012D7">
      <w:rPr>
        <w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cs="Tahoma"/>
        <w:i/>
        <w:iCs/>
        <w:sz w:val="20"/>
        <w:szCs w:val="20"/>
        <w:lang w:val="en-US"/>
      </w:rPr>
      <w:t>May 15, 2010</w:t>
    </w:r>
  </w:smartTag>
</w:p>
"@
$dateSearch = "Date: " + " " + [char]0x2013 + " "
$pDate = Find-ParagraphByText($dateSearch)
$pDate.Range.InsertXML($dateXml)

# ---------------------------------------------------------------------------
# 2. "2009-08 - Review by Ramses " -> "2009-08 - Review by Brother "
# ---------------------------------------------------------------------------
$reviewXml = @"
<w:p $WNS>
  <w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr>
  <w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">2009-08 - Review by </w:t></w:r>
  <w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Brother</w:t></w:r>
  <w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>
</w:p>
"@
$pReview = Find-ParagraphByText("2009-08 - Review by Ramses ")
$pReview.Range.InsertXML($reviewXml)

# ---------------------------------------------------------------------------
# 3. "2010-05 - Events" -> "2010-05 - Events" (re-run-split) followed by a
#    brand-new "2010-05 - Inheritance" paragraph
# ---------------------------------------------------------------------------
$eventsXml = @"
<w:p $WNS>
  <w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr>
  <w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">2010-05 </w:t></w:r>
  <w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>-</w:t></w:r>
  <w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> Events</w:t></w:r>
</w:p>
<w:p $WNS>
  <w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr>
  <w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">2010-05 - </w:t></w:r>
  <w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Inheritance</w:t></w:r>
</w:p>
"@
$pEvents = Find-ParagraphByText("2010-05 - Events")
$pEvents.Range.InsertXML($eventsXml)

# ---------------------------------------------------------------------------
# 4. "What is nice to realize ... actually easily worked out:" ->
#    "... actually more easily worked out than expected:" + a new blank
#    paragraph after it
# ---------------------------------------------------------------------------
$niceXml = @"
<w:p $WNS>
  <w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr>
  <w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">What is nice to realize is that one &#8216;hard topic&#8217; was actually </w:t></w:r>
  <w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">more </w:t></w:r>
  <w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>easily worked out</w:t></w:r>
  <w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> than expected</w:t></w:r>
  <w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>:</w:t></w:r>
</w:p>
<w:p $WNS>
  <w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr>
</w:p>
"@
$pNice = Find-ParagraphByText("What is nice to realize is that one " + [char]0x2018 + "hard topic" + [char]0x2019 + " was actually easily worked out:")
$pNice.Range.InsertXML($niceXml)

# ---------------------------------------------------------------------------
# 5. "- Reliability of interfaces in dual relations" ->
#    "- Reliability of interfaces in bi-directional relationships" followed
#    by a new blank (indented) paragraph
# ---------------------------------------------------------------------------
$reliabilityXml = @"
<w:p $WNS>
  <w:pPr><w:ind w:left="284"/><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr>
  <w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">- Reliability of interfaces in </w:t></w:r>
  <w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>bi-directional</w:t></w:r>
  <w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> relations</w:t></w:r>
  <w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>hips</w:t></w:r>
</w:p>
<w:p $WNS>
  <w:pPr><w:ind w:left="284"/><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr>
</w:p>
"@
$pReliability = Find-ParagraphByText("- Reliability of interfaces in dual relations")
$pReliability.Range.InsertXML($reliabilityXml)

# ---------------------------------------------------------------------------
# 6. "This hard issue was gracefully and easily solved during the
#    Interfaces project." ->
#    "This was considered a difficult topic, and adequately solved during
#    the Interfaces project."
# ---------------------------------------------------------------------------
$solvedXml = @"
<w:p $WNS>
  <w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr>
  <w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">This was </w:t></w:r>
  <w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">considered a </w:t></w:r>
  <w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>difficult</w:t></w:r>
  <w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> topic</w:t></w:r>
  <w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">, and </w:t></w:r>
  <w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">adequately </w:t></w:r>
  <w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>solved during the Interfaces project.</w:t></w:r>
</w:p>
"@
$pSolved = Find-ParagraphByText("This hard issue was gracefully and easily solved during the Interfaces project.")
$pSolved.Range.InsertXML($solvedXml)

# ---------------------------------------------------------------------------
# 7. "So when you focus on the easy things, it gives the hard things time
#    to settle and by the time you get to them the penny may already have
#    dropped." ->
#    "So perhaps first trying to focus on the 'easier' topics, might give
#    the hard topics time to settle and by the time you get to them, the
#    idea may be clearer."
#    (the _GoBack bookmark sits between "them," and " the idea may ...")
# ---------------------------------------------------------------------------
$pennyXml = @"
<w:p $WNS>
  <w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr>
  <w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">So </w:t></w:r>
  <w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">perhaps first </w:t></w:r>
  <w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">trying </w:t></w:r>
  <w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">to </w:t></w:r>
  <w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">focus on the </w:t></w:r>
  <w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>'easier' topics</w:t></w:r>
  <w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">, </w:t></w:r>
  <w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">might </w:t></w:r>
  <w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">give the hard </w:t></w:r>
  <w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>topics</w:t></w:r>
  <w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> time to settle and by the time you get to them</w:t></w:r>
  <w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>,</w:t></w:r>
  <w:bookmarkStart w:id="0" w:name="_GoBack"/>
  <w:bookmarkEnd w:id="0"/>
  <w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> the idea may </w:t></w:r>
  <w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">be </w:t></w:r>
  <w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>clearer</w:t></w:r>
  <w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>.</w:t></w:r>
</w:p>
"@
$pPenny = Find-ParagraphByText("So when you focus on the easy things, it gives the hard things time to settle and by the time you get to them the penny may already have dropped.")
$pPenny.Range.InsertXML($pennyXml)

Write-Output "Done. Paragraph count now: $($d.Paragraphs.Count)"
